# Updates cryptos price (D) and 1h-volume-change (E) columns to match the
# refreshed GitHub Actions data snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Preserve the cell's existing style while forcing the assigned value to
    # be stored as text (not auto-coerced to a number), matching the source
    # workbook where every Price cell is an inline string.
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $savedStyle
}

$ws.Range("D2").Value = "43.417.21"
$ws.Range("E2").Value = "  +2.78%  "

$ws.Range("D3").Value = "2.308.60"
$ws.Range("E3").Value = "  +1.73%  "

$ws.Range("E4").Value = "  +0.03%  "

Set-TextValue $ws.Range("D5") "311.27"
$ws.Range("E5").Value = "  +1.48%  "

Set-TextValue $ws.Range("D6") "102.75"
$ws.Range("E6").Value = "  +5.52%  "

$ws.Range("E7").Value = "  +1.29%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +7.74%  "

Set-TextValue $ws.Range("D10") "35.85"
$ws.Range("E10").Value = "  +1.77%  "

Set-TextValue $ws.Range("D11") "0.0813"
$ws.Range("E11").Value = "  +2.75%  "

$ws.Range("E12").Value = "  -0.94%  "

Set-TextValue $ws.Range("D13") "6.99"
$ws.Range("E13").Value = "  +1.69%  "

$ws.Range("D14").Value = "2.665.76"
$ws.Range("E14").Value = "  +1.79%  "

Set-TextValue $ws.Range("D15") "15.03"
$ws.Range("E15").Value = "  +2.23%  "

$ws.Range("D16").Value = "2.314.84"
$ws.Range("E16").Value = "  +2.15%  "

Set-TextValue $ws.Range("D17") "0.807"
$ws.Range("E17").Value = "  +1.89%  "

$ws.Range("D18").Value = "43.313.38"
$ws.Range("E18").Value = "  +2.91%  "

Set-TextValue $ws.Range("D19") "12.33"
$ws.Range("E19").Value = "  +0.42%  "

$ws.Range("E20").Value = "  +3.28%  "

Set-TextValue $ws.Range("D21") "6.18"
$ws.Range("E21").Value = "  +2.69%  "

Set-TextValue $ws.Range("D22") "68.13"
$ws.Range("E22").Value = "  +0.44%  "

Set-TextValue $ws.Range("D23") "241.59"
$ws.Range("E23").Value = "  +1.79%  "

$ws.Range("E24").Value = "  +1.40%  "

Set-TextValue $ws.Range("D25") "2.01"
$ws.Range("E25").Value = "  +2.22%  "

$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("E27").Value = "  -1.85%  "

Set-TextValue $ws.Range("D28") "24.73"
$ws.Range("E28").Value = "  +4.91%  "

$ws.Range("E29").Value = "  +8.50%  "

Set-TextValue $ws.Range("D30") "36.82"
$ws.Range("E30").Value = "  -2.75%  "

Set-TextValue $ws.Range("D31") "9.64"
$ws.Range("E31").Value = "  +0.47%  "

Set-TextValue $ws.Range("D32") "168.43"
$ws.Range("E32").Value = "  +3.64%  "

Set-TextValue $ws.Range("D33") "5.29"
$ws.Range("E33").Value = "  +0.63%  "

$ws.Range("E34").Value = "  +0.04%  "

Set-TextValue $ws.Range("D35") "2.52"
$ws.Range("E35").Value = "  +6.17%  "

$ws.Range("E36").Value = "  +0.82%  "

Set-TextValue $ws.Range("D37") "17.76"
$ws.Range("E37").Value = "  +0.06%  "

$ws.Range("E38").Value = "  -2.74%  "

Set-TextValue $ws.Range("D39") "1.91"
$ws.Range("E39").Value = "  +4.68%  "

$ws.Range("E40").Value = "  +1.93%  "

$ws.Range("E41").Value = "  +1.35%  "

$ws.Range("E42").Value = "  +6.31%  "

Set-TextValue $ws.Range("D43") "2.31"
$ws.Range("E43").Value = "  -1.47%  "

Set-TextValue $ws.Range("D44") "19.58"
$ws.Range("E44").Value = "  +2.52%  "

$ws.Range("E45").Value = "  +2.93%  "

$ws.Range("D46").Value = "1.968.70"
$ws.Range("E46").Value = "  +0.78%  "

Set-TextValue $ws.Range("D47") "3.00"
$ws.Range("E47").Value = "  +2.55%  "

Set-TextValue $ws.Range("D48") "9.87"
$ws.Range("E48").Value = "  -1.48%  "

Set-TextValue $ws.Range("D49") "55.51"
$ws.Range("E49").Value = "  +2.68%  "

Set-TextValue $ws.Range("D50") "1.58"
$ws.Range("E50").Value = "  +7.35%  "

Set-TextValue $ws.Range("D51") "2.90"
$ws.Range("E51").Value = "  +0.97%  "
